# "modified the model and milp"
#
# Renames the four variable-name headers on Sheet1 (row 1, columns C:F)
# from the old P_EL / H_TL / C_TL / P_RES naming scheme to the new
# L_Power / L_Heat / L_Cool / P_PV naming scheme, and moves the active
# selection on Sheet1 from H9 to G3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the model/MILP variable headers (shared-string text) in row 1.
$ws.Range("C1").Value = "L_Power"
$ws.Range("D1").Value = "L_Heat"
$ws.Range("E1").Value = "L_Cool"
$ws.Range("F1").Value = "P_PV"

# Move the active selection to match the updated working cell.
$ws.Activate() | Out-Null
$ws.Range("G3").Select() | Out-Null
